$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.9738157391548157
$ws.Range("B1").Value = 1.625557899475098
$ws.Range("C1").Value = 6.72291898727417
$ws.Range("D1").Value = 2.679537773132324
$ws.Range("E1").Value = 1.50115966796875
